# Added code for return flights testing.
#
# Summary of changes applied to DataSheets/Data_sheets_expedia.xlsx:
#  1. Rename sheet "OneWayFlightsPosDefaultDate" -> "FlightsPosDefaultDate"
#     and update its tag-markers (A1 / D6) from "OneWayFlightsPos" to the
#     new tag name "FlightsDefaultDtPos". Also move its selection to I9:J10.
#  2. Delete the sheet "OneWayFlightsTravellersAge" entirely (its
#     functionality/data now lives on the PosAllOptions sheet which already
#     supports the return-flight / traveller-age columns).
#  3. Replace the rich (mixed-bold) instructions text on the
#     "Important Information" sheet with a simplified, plain-text
#     instruction paragraph.
#  4. Nudge the view of "OneWayFlightsPosAllOptions" so column D is the
#     left-most visible column (extra columns were added for the new
#     return-flight related fields), keeping K7 selected/active.

$excel.DisplayAlerts = $false | Out-Null

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) OneWayFlightsPosDefaultDate -> FlightsPosDefaultDate
# ---------------------------------------------------------------------
$wsDefaultDate = $wb.Worksheets.Item("OneWayFlightsPosDefaultDate")
$wsDefaultDate.Name = "FlightsPosDefaultDate"

$wsDefaultDate.Range("A1").Value = "FlightsDefaultDtPos"
$wsDefaultDate.Range("D6").Value = "FlightsDefaultDtPos"

$wsDefaultDate.Select() | Out-Null
$wsDefaultDate.Range("I9:J10").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the OneWayFlightsTravellersAge sheet
# ---------------------------------------------------------------------
$wsTravellersAge = $wb.Worksheets.Item("OneWayFlightsTravellersAge")
$wsTravellersAge.Delete() | Out-Null

# ---------------------------------------------------------------------
# 3) Important Information -> simplified plain-text instructions
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Important Information")
$newInstructions = "All sheets with the word 'Pos' within their tab name indicate test data for positive test cases." + [char]10 + "Look out for specific instruction in each sheet."
$wsInfo.Range("A1").Value = $newInstructions

# ---------------------------------------------------------------------
# 4) OneWayFlightsPosAllOptions view tweak - scroll so column D is first,
#    keep the previously active selection/tab.
# ---------------------------------------------------------------------
$wsAllOptions = $wb.Worksheets.Item("OneWayFlightsPosAllOptions")
$wsAllOptions.Select() | Out-Null
$wsAllOptions.Range("K7").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4

Write-Host "Edit complete"
